# Updates cryptos list with latest prices / volume(1h) data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.364.79"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "3.473.98"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D5").Value = "'593.80"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'179.27"
$ws.Range("E6").Value = "  +4.64%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.474.35"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  +5.86%  "
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "'0.432"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "4.075.59"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "'32.21"
$ws.Range("E14").Value = "  +12.19%  "
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "67.368.66"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "3.472.48"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "'14.30"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("D21").Value = "'389.85"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'72.83"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'5.73"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "'0.534"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "'10.34"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'6.19"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'23.51"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "'7.38"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").Value = "'163.67"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'0.870"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "'1.88"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = "  +7.52%  "
$ws.Range("D42").Value = "'6.88"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "2.833.12"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "'0.0723"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").Value = "'26.58"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").Value = "'0.0299"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "'336.24"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -1.83%  "
